# Update "horarios" workbook with the latest scrape for Línea 141.
# New scrape time: 01:56:45 (was 01:16:06)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "LP1912": update existing rows 6-7, append a new row 8, and
# refresh the header (timestamp + row count).
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 01:56:45"
$ws1.Range("A3").Value = "Total filas: 3"

# Row 6: 14_ABASTO — arrival time unchanged, minutes refreshed
$ws1.Range("A6").Value = "01:56:45"
$ws1.Range("B6").Value = "01:58"
$ws1.Range("C6").Value = "14_ABASTO"
$ws1.Range("D6").Value = 2
$ws1.Range("E6").Value = "LP1912"

# Row 7: 215_ALUAR — refreshed arrival time + minutes
$ws1.Range("A7").Value = "01:56:45"
$ws1.Range("B7").Value = "03:04"
$ws1.Range("C7").Value = "215_ALUAR"
$ws1.Range("D7").Value = 68
$ws1.Range("E7").Value = "LP1912"

# Row 8 (new): 14_ABASTO second arrival
$ws1.Range("A8").Value = "01:56:45"
$ws1.Range("B8").Value = "03:48"
$ws1.Range("C8").Value = "14_ABASTO"
$ws1.Range("D8").Value = 112
$ws1.Range("E8").Value = "LP1912"

# ---------------------------------------------------------------------
# Sheet "LP1912-215": refresh header timestamp and the single data row.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 01:56:45"

$ws2.Range("A6").Value = "01:56:45"
$ws2.Range("B6").Value = "03:04"
$ws2.Range("C6").Value = "215_ALUAR"
$ws2.Range("D6").Value = 68
$ws2.Range("E6").Value = "LP1912"

# ---------------------------------------------------------------------
# Sheet "6203-6173": no data rows, only the timestamp header changes.
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 01:56:45"
